$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 2
$ws.Range("B6").Value = 2
$ws.Range("B7").Value = 3
$ws.Range("B8").Value = 2
$ws.Range("B9").Value = 2
$ws.Range("B10").Value = 2
$ws.Range("B11").Value = 5
$ws.Range("B13").Value = 5
$ws.Range("B14").Value = 5
$ws.Range("B15").Value = 2
$ws.Range("B16").Value = 2
$ws.Range("B17").Value = 4
$ws.Range("B19").Value = 5
$ws.Range("B20").Value = 4
$ws.Range("B21").Value = 5
$ws.Range("B22").Value = 4
$ws.Range("B23").Value = 1

$ws.Range("B24").Select()
